# Refresh the crypto price ("D") and 1h-volume-change ("E") columns for
# rows 2-51 with the latest scraped figures (GitHub Actions data refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings such as "58.252.30" or "135.00" as plain text
# cells (no numeric NumberFormat applied, t="inlineStr" in the OOXML). Writing
# such a string straight into .Value lets Excel's smart-parser coerce a
# numeric-looking value (e.g. "135.00") into a real Number, which silently
# drops the trailing zero and flips the cell's stored type. Prefixing the
# value with a literal leading apostrophe forces Excel to keep it as text
# (exactly like typing  '135.00  into the cell by hand); the Style reset that
# follows clears the resulting quote-prefix formatting so the cell keeps the
# sheet's default (unstyled) look, matching every other data cell.
function Set-TextCell($addr, $val) {
    $ws.Range($addr).Value = "'" + $val
    $ws.Range($addr).Style = "Normal"
}

Set-TextCell "D2" '58.252.30'
$ws.Range("E2").Value = '  +1.77%  '
Set-TextCell "D3" '2.345.30'
$ws.Range("E3").Value = '  +0.68%  '
Set-TextCell "D5" '542.16'
$ws.Range("E5").Value = '  +2.09%  '
Set-TextCell "D6" '135.00'
$ws.Range("E6").Value = '  +1.84%  '
$ws.Range("E7").Value = '  +0.63%  '
Set-TextCell "D8" '0.561'
$ws.Range("E8").Value = '  +4.88%  '
$ws.Range("E9").Value = '  +0.41%  '
Set-TextCell "D10" '5.67'
$ws.Range("E10").Value = '  +6.75%  '
$ws.Range("E11").Value = '  -0.65%  '
Set-TextCell "D12" '0.356'
$ws.Range("E12").Value = '  +3.56%  '
Set-TextCell "D13" '23.83'
$ws.Range("E13").Value = '  +1.05%  '
Set-TextCell "D14" '2.761.32'
$ws.Range("E14").Value = '  +0.43%  '
Set-TextCell "D15" '58.148.13'
$ws.Range("E15").Value = '  +1.56%  '
Set-TextCell "D16" '0.0000134'
$ws.Range("E16").Value = '  +0.23%  '
Set-TextCell "D17" '2.318.83'
$ws.Range("E17").Value = '  -1.29%  '
Set-TextCell "D18" '10.73'
$ws.Range("E18").Value = '  +2.64%  '
Set-TextCell "D19" '333.24'
$ws.Range("E19").Value = '  -1.36%  '
Set-TextCell "D20" '4.26'
$ws.Range("E20").Value = '  +1.95%  '
Set-TextCell "D21" '6.66'
$ws.Range("E21").Value = '  -3.78%  '
$ws.Range("E22").Value = '  +0.06%  '
$ws.Range("E23").Value = '  +0.05%  '
Set-TextCell "D24" '62.76'
$ws.Range("E24").Value = '  +1.59%  '
Set-TextCell "D26" '8.54'
$ws.Range("E26").Value = '  -3.95%  '
Set-TextCell "D27" '0.999'
$ws.Range("E27").Value = '  +0.90%  '
Set-TextCell "D28" '1.42'
$ws.Range("E28").Value = '  +5.44%  '
$ws.Range("E29").Value = '  +2.10%  '
Set-TextCell "D30" '170.28'
$ws.Range("E30").Value = '  +0.19%  '
Set-TextCell "D31" '0.0₃0737'
$ws.Range("E31").Value = '  +1.32%  '
Set-TextCell "D32" '6.11'
$ws.Range("E32").Value = '  -0.38%  '
$ws.Range("E33").Value = '  +12.60%  '
$ws.Range("E34").Value = '  -0.53%  '
$ws.Range("E35").Value = '  +0.03%  '
Set-TextCell "D36" '4.25'
$ws.Range("E36").Value = '  +5.46%  '
Set-TextCell "D37" '0.999'
$ws.Range("E37").Value = '  +0.91%  '
Set-TextCell "D38" '1.25'
$ws.Range("E38").Value = '  -1.70%  '
$ws.Range("E39").Value = '  +3.29%  '
Set-TextCell "D40" '39.11'
$ws.Range("E40").Value = '  +0.35%  '
Set-TextCell "D41" '142.15'
$ws.Range("E41").Value = '  -4.09%  '
$ws.Range("E42").Value = '  +1.57%  '
$ws.Range("E43").Value = '  -0.37%  '
Set-TextCell "D44" '287.22'
$ws.Range("E44").Value = '  +0.40%  '
$ws.Range("E45").Value = '  +0.53%  '
Set-TextCell "D46" '19.20'
$ws.Range("E46").Value = '  +2.07%  '
$ws.Range("E47").Value = '  -0.03%  '
Set-TextCell "D48" '0.565'
$ws.Range("E48").Value = '  +0.52%  '
$ws.Range("E49").Value = '  +0.46%  '
Set-TextCell "D50" '0.381'
$ws.Range("E50").Value = '  +0.92%  '
$ws.Range("E51").Value = '  +0.46%  '
